# Apply the daily cryptos-list refresh (prices + 1h volume%, and a few
# rank swaps / a dropped-in replacement coin) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.264.70"
$ws.Range("E2").Value = "  +2.72%  "

$ws.Range("D3").Value = "2.070.53"
$ws.Range("E3").Value = "  +2.87%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'251.56"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'0.648"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").Value = "'64.24"
$ws.Range("E7").Value = "  +1.57%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.394"
$ws.Range("E9").Value = "  +5.61%  "

$ws.Range("D10").Value = "'59.49"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  +8.56%  "

$ws.Range("D12").Value = "'0.104"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "'0.920"
$ws.Range("E13").Value = "  -2.95%  "

$ws.Range("D14").Value = "'23.20"
$ws.Range("E14").Value = "  +17.36%  "

$ws.Range("D15").Value = "'14.71"
$ws.Range("E15").Value = "  -2.54%  "

$ws.Range("D16").Value = "2.375.80"
$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("D17").Value = "'5.62"
$ws.Range("E17").Value = "  +2.98%  "

$ws.Range("D18").Value = "2.078.31"
$ws.Range("E18").Value = "  +3.09%  "

$ws.Range("D19").Value = "37.293.18"
$ws.Range("E19").Value = "  +2.88%  "

$ws.Range("D20").Value = "'73.43"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("D21").Value = "0.0₃0898"
$ws.Range("E21").Value = "  +4.44%  "

$ws.Range("D22").Value = "'5.49"
$ws.Range("E22").Value = "  +3.72%  "

$ws.Range("D23").Value = "'238.55"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  -3.09%  "

$ws.Range("E26").Value = "  +2.54%  "

$ws.Range("D27").Value = "'9.96"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D28").Value = "'20.64"
$ws.Range("E28").Value = "  +4.85%  "

$ws.Range("D29").Value = "'160.97"
$ws.Range("E29").Value = "  -3.19%  "

$ws.Range("E30").Value = "  +28.74%  "

$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("D32").Value = "'5.13"
$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").Value = "'1.19"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("D34").Value = "'0.0626"
$ws.Range("E34").Value = "  +2.80%  "

$ws.Range("D35").Value = "'4.63"
$ws.Range("E35").Value = "  +2.90%  "

$ws.Range("D36").Value = "'2.52"
$ws.Range("E36").Value = "  +1.64%  "

$ws.Range("D37").Value = "'6.40"
$ws.Range("E37").Value = "  +10.46%  "

$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.84"
$ws.Range("E39").Value = "  +1.64%  "

$ws.Range("D40").Value = "'3.02"
$ws.Range("E40").Value = "  +28.80%  "

$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.102"
$ws.Range("E41").Value = "  +5.79%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.28"
$ws.Range("E42").Value = "  +2.73%  "

$ws.Range("D44").Value = "'1.17"
$ws.Range("E44").Value = "  +3.93%  "

$ws.Range("D45").Value = "'17.21"
$ws.Range("E45").Value = "  +0.84%  "

$ws.Range("D46").Value = "'0.0218"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").Value = "'7.96"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("D48").Value = "'95.44"
$ws.Range("E48").Value = "  +0.56%  "

$ws.Range("D49").Value = "1.397.82"
$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("D50").Value = "'2.93"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.257.23"
$ws.Range("E51").Value = "  +2.96%  "
